# Append the new match row (row 99) to the Ecuador Liga Pro 2023 sheet,
# mirroring the formatting/style already used by the preceding data row (98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 99
$prevRow = 98

# Column D ("temporada") holds the numeric-looking text "2023" as a real
# string in the source data (same as every other row). Force text storage
# via a temporary "@" number format before assigning the value so COM
# doesn't silently coerce it to a number; the format-copy step below then
# restores the normal (General) formatting from the previous row.
$ws.Cells.Item($newRow, 4).NumberFormat = "@"

# Fill in the values for the new row.
$ws.Cells.Item($newRow, 1).Value = 98
$ws.Cells.Item($newRow, 2).Value = "ecuador"
$ws.Cells.Item($newRow, 3).Value = "liga-pro"
$ws.Cells.Item($newRow, 4).Value = "2023"
$ws.Cells.Item($newRow, 5).Value = 45241.79166666666
$ws.Cells.Item($newRow, 6).Value = "U. Catolica"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Cumbaya"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 1.32
$ws.Cells.Item($newRow, 11).Value = "06/11/2023 01:12"
$ws.Cells.Item($newRow, 12).Value = 1.37
$ws.Cells.Item($newRow, 13).Value = "11/11/2023 18:56"
$ws.Cells.Item($newRow, 14).Value = 4.96
$ws.Cells.Item($newRow, 15).Value = "06/11/2023 01:12"
$ws.Cells.Item($newRow, 16).Value = 4.9
$ws.Cells.Item($newRow, 17).Value = "11/11/2023 18:56"
$ws.Cells.Item($newRow, 18).Value = 8.220000000000001
$ws.Cells.Item($newRow, 19).Value = "06/11/2023 01:12"
$ws.Cells.Item($newRow, 20).Value = 8.48
$ws.Cells.Item($newRow, 21).Value = "11/11/2023 18:56"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/ecuador/liga-pro/u-catolica-cumbaya/tGaAW5ZG/"

# Copy the formatting (styles) of the last existing data row onto the new
# row so the new cells pick up the same style indices (bold/border for
# column A, date-time number format for column E, General elsewhere) as
# the rest of the table. Done last so it doesn't get clobbered by, and
# doesn't clobber, the values set above.
$ws.Range("A$prevRow`:V$prevRow").Copy() | Out-Null
$ws.Range("A$newRow`:V$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
